# Updates cryptos list prices and 1h volume percentages (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.767.92'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.145.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.14%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.43%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.143.89'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.10'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.498'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('E13').Value = '  +0.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.659.62'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.914.84'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.145.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.35%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '500.50'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('E22').Value = '  -2.60%  '
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.16%  '
$ws.Range('E26').Value = '  -0.27%  '
$ws.Range('E27').Value = '  -1.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.81'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.82'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.04%  '
$ws.Range('E35').Value = '  -2.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '54.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0890'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '475.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('E40').Value = '  -2.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.61'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.000.85'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.46%  '
$ws.Range('E43').Value = '  -3.66%  '
$ws.Range('E44').Value = '  -3.40%  '
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₃0578'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.14'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.96%  '
